$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 3).Value = 56232
$ws.Cells.Item(2, 4).Value = 114710489
$ws.Cells.Item(3, 3).Value = 136430
$ws.Cells.Item(3, 4).Value = 320513794
$ws.Cells.Item(4, 3).Value = 49522
$ws.Cells.Item(4, 4).Value = 143820446
$ws.Cells.Item(5, 3).Value = 15682
$ws.Cells.Item(5, 4).Value = 53020452
$ws.Cells.Item(6, 3).Value = 5770
$ws.Cells.Item(6, 4).Value = 26271871
$ws.Cells.Item(7, 3).Value = 1147
$ws.Cells.Item(7, 4).Value = 6688434
$ws.Cells.Item(8, 3).Value = 70
$ws.Cells.Item(8, 4).Value = 471481
$ws.Cells.Item(12, 3).Value = 58502
$ws.Cells.Item(12, 4).Value = 93094711
$ws.Cells.Item(13, 3).Value = 14269
$ws.Cells.Item(13, 4).Value = 28744342
$ws.Cells.Item(14, 3).Value = 38393
$ws.Cells.Item(14, 4).Value = 88628732
$ws.Cells.Item(15, 3).Value = 12753
$ws.Cells.Item(15, 4).Value = 35299250
$ws.Cells.Item(16, 3).Value = 3671
$ws.Cells.Item(16, 4).Value = 11262269
$ws.Cells.Item(17, 3).Value = 1199
$ws.Cells.Item(17, 4).Value = 5157416
$ws.Cells.Item(18, 3).Value = 268
$ws.Cells.Item(18, 4).Value = 1449487
$ws.Cells.Item(20, 3).Value = 14333
$ws.Cells.Item(20, 4).Value = 22222419
$ws.Cells.Item(21, 3).Value = 20047
$ws.Cells.Item(21, 4).Value = 42393023
$ws.Cells.Item(22, 3).Value = 47515
$ws.Cells.Item(22, 4).Value = 114256875
$ws.Cells.Item(23, 3).Value = 16456
$ws.Cells.Item(23, 4).Value = 47445108
$ws.Cells.Item(24, 3).Value = 4869
$ws.Cells.Item(24, 4).Value = 15736377
$ws.Cells.Item(25, 3).Value = 1563
$ws.Cells.Item(25, 4).Value = 6523892
$ws.Cells.Item(28, 3).Value = 16037
$ws.Cells.Item(28, 4).Value = 24728200
$ws.Cells.Item(29, 3).Value = 11342
$ws.Cells.Item(29, 4).Value = 23167236
$ws.Cells.Item(30, 3).Value = 32840
$ws.Cells.Item(30, 4).Value = 74200150
$ws.Cells.Item(31, 3).Value = 11852
$ws.Cells.Item(31, 4).Value = 32067377
$ws.Cells.Item(32, 3).Value = 3258
$ws.Cells.Item(32, 4).Value = 9755112
$ws.Cells.Item(33, 3).Value = 1023
$ws.Cells.Item(33, 4).Value = 4334545
$ws.Cells.Item(34, 3).Value = 205
$ws.Cells.Item(34, 4).Value = 961324
$ws.Cells.Item(36, 3).Value = 11602
$ws.Cells.Item(36, 4).Value = 18023081
$ws.Cells.Item(37, 3).Value = 5016
$ws.Cells.Item(37, 4).Value = 11038077
$ws.Cells.Item(38, 3).Value = 11784
$ws.Cells.Item(38, 4).Value = 27393981
$ws.Cells.Item(39, 3).Value = 4874
$ws.Cells.Item(39, 4).Value = 13872526
$ws.Cells.Item(40, 3).Value = 1355
$ws.Cells.Item(40, 4).Value = 4427400
$ws.Cells.Item(41, 3).Value = 433
$ws.Cells.Item(41, 4).Value = 2076184
$ws.Cells.Item(42, 3).Value = 52
$ws.Cells.Item(42, 4).Value = 329772
$ws.Cells.Item(44, 3).Value = 3570
$ws.Cells.Item(44, 4).Value = 5496539
$ws.Cells.Item(45, 3).Value = 25699
$ws.Cells.Item(45, 4).Value = 53145514
$ws.Cells.Item(46, 3).Value = 76274
$ws.Cells.Item(46, 4).Value = 179920123
$ws.Cells.Item(47, 3).Value = 29347
$ws.Cells.Item(47, 4).Value = 81916791
$ws.Cells.Item(48, 3).Value = 9550
$ws.Cells.Item(48, 4).Value = 29099305
$ws.Cells.Item(49, 3).Value = 3287
$ws.Cells.Item(49, 4).Value = 13214539
$ws.Cells.Item(50, 3).Value = 569
$ws.Cells.Item(50, 4).Value = 3184876
$ws.Cells.Item(53, 3).Value = 26346
$ws.Cells.Item(53, 4).Value = 48368713
$ws.Cells.Item(54, 3).Value = 2707
$ws.Cells.Item(54, 4).Value = 4387992
$ws.Cells.Item(55, 3).Value = 9040
$ws.Cells.Item(55, 4).Value = 14985554
$ws.Cells.Item(56, 3).Value = 3032
$ws.Cells.Item(56, 4).Value = 5346977
$ws.Cells.Item(57, 3).Value = 994
$ws.Cells.Item(57, 4).Value = 1949513
$ws.Cells.Item(61, 3).Value = 9243
$ws.Cells.Item(61, 4).Value = 13747465
$ws.Cells.Item(62, 3).Value = 1825
$ws.Cells.Item(62, 4).Value = 3980094
$ws.Cells.Item(63, 3).Value = 4316
$ws.Cells.Item(63, 4).Value = 9361700
$ws.Cells.Item(64, 3).Value = 1727
$ws.Cells.Item(64, 4).Value = 3886459
$ws.Cells.Item(66, 3).Value = 210
$ws.Cells.Item(66, 4).Value = 471487
$ws.Cells.Item(68, 3).Value = 2828
$ws.Cells.Item(68, 4).Value = 5597206
$ws.Cells.Item(69, 3).Value = 22889
$ws.Cells.Item(69, 4).Value = 45278224
$ws.Cells.Item(70, 3).Value = 66425
$ws.Cells.Item(70, 4).Value = 151647727
$ws.Cells.Item(71, 3).Value = 24289
$ws.Cells.Item(71, 4).Value = 67498653
$ws.Cells.Item(72, 3).Value = 7601
$ws.Cells.Item(72, 4).Value = 23010657
$ws.Cells.Item(73, 3).Value = 2465
$ws.Cells.Item(73, 4).Value = 9967043
$ws.Cells.Item(74, 3).Value = 496
$ws.Cells.Item(74, 4).Value = 2742238
$ws.Cells.Item(75, 3).Value = 27
$ws.Cells.Item(75, 4).Value = 88619
$ws.Cells.Item(78, 3).Value = 21234
$ws.Cells.Item(78, 4).Value = 32615041
$ws.Cells.Item(79, 3).Value = 83469
$ws.Cells.Item(79, 4).Value = 171523501
$ws.Cells.Item(80, 3).Value = 226810
$ws.Cells.Item(80, 4).Value = 511707455
$ws.Cells.Item(81, 3).Value = 102420
$ws.Cells.Item(81, 4).Value = 287462800
$ws.Cells.Item(82, 3).Value = 37346
$ws.Cells.Item(82, 4).Value = 125670917
$ws.Cells.Item(83, 3).Value = 13776
$ws.Cells.Item(83, 4).Value = 62319089
$ws.Cells.Item(84, 3).Value = 2697
$ws.Cells.Item(84, 4).Value = 17333162
$ws.Cells.Item(90, 3).Value = 79255
$ws.Cells.Item(90, 4).Value = 125991855
$ws.Cells.Item(91, 3).Value = 5634
$ws.Cells.Item(91, 4).Value = 8785893
$ws.Cells.Item(92, 3).Value = 13575
$ws.Cells.Item(92, 4).Value = 21494096
$ws.Cells.Item(93, 3).Value = 4357
$ws.Cells.Item(93, 4).Value = 7076451
$ws.Cells.Item(95, 3).Value = 419
$ws.Cells.Item(95, 4).Value = 927812
$ws.Cells.Item(98, 3).Value = 6362
$ws.Cells.Item(98, 4).Value = 8757623
$ws.Cells.Item(99, 3).Value = 2310
$ws.Cells.Item(99, 4).Value = 4277260
$ws.Cells.Item(100, 3).Value = 7328
$ws.Cells.Item(100, 4).Value = 14577009
$ws.Cells.Item(101, 3).Value = 2629
$ws.Cells.Item(101, 4).Value = 5951843
$ws.Cells.Item(103, 3).Value = 339
$ws.Cells.Item(103, 4).Value = 1080276
$ws.Cells.Item(104, 3).Value = 68
$ws.Cells.Item(104, 4).Value = 304891
$ws.Cells.Item(106, 3).Value = 4913
$ws.Cells.Item(106, 4).Value = 7168362
$ws.Cells.Item(107, 3).Value = 1038
$ws.Cells.Item(107, 4).Value = 2358697
$ws.Cells.Item(113, 3).Value = 16402
$ws.Cells.Item(113, 4).Value = 34478606
$ws.Cells.Item(114, 3).Value = 43177
$ws.Cells.Item(114, 4).Value = 101175041
$ws.Cells.Item(115, 3).Value = 15318
$ws.Cells.Item(115, 4).Value = 42874672
$ws.Cells.Item(116, 3).Value = 4732
$ws.Cells.Item(116, 4).Value = 14986972
$ws.Cells.Item(117, 3).Value = 1466
$ws.Cells.Item(117, 4).Value = 6207982
$ws.Cells.Item(118, 3).Value = 304
$ws.Cells.Item(118, 4).Value = 1725408
$ws.Cells.Item(122, 3).Value = 13527
$ws.Cells.Item(122, 4).Value = 20745951
$ws.Cells.Item(123, 3).Value = 44704
$ws.Cells.Item(123, 4).Value = 90644083
$ws.Cells.Item(124, 3).Value = 95673
$ws.Cells.Item(124, 4).Value = 215440980
$ws.Cells.Item(125, 3).Value = 32350
$ws.Cells.Item(125, 4).Value = 87050498
$ws.Cells.Item(126, 3).Value = 10137
$ws.Cells.Item(126, 4).Value = 30842374
$ws.Cells.Item(127, 3).Value = 3211
$ws.Cells.Item(127, 4).Value = 13185501
$ws.Cells.Item(128, 3).Value = 656
$ws.Cells.Item(128, 4).Value = 3561327
$ws.Cells.Item(132, 3).Value = 35060
$ws.Cells.Item(132, 4).Value = 53649627
$ws.Cells.Item(133, 3).Value = 53793
$ws.Cells.Item(133, 4).Value = 110422672
$ws.Cells.Item(134, 3).Value = 112625
$ws.Cells.Item(134, 4).Value = 251657072
$ws.Cells.Item(135, 3).Value = 36486
$ws.Cells.Item(135, 4).Value = 100462357
$ws.Cells.Item(136, 3).Value = 10794
$ws.Cells.Item(136, 4).Value = 33326206
$ws.Cells.Item(137, 3).Value = 3397
$ws.Cells.Item(137, 4).Value = 14012897
$ws.Cells.Item(138, 3).Value = 546
$ws.Cells.Item(138, 4).Value = 2996655
$ws.Cells.Item(139, 3).Value = 47
$ws.Cells.Item(139, 4).Value = 232914
$ws.Cells.Item(142, 3).Value = 43789
$ws.Cells.Item(142, 4).Value = 65645814
$ws.Cells.Item(143, 3).Value = 19651
$ws.Cells.Item(143, 4).Value = 40380384
$ws.Cells.Item(144, 3).Value = 47766
$ws.Cells.Item(144, 4).Value = 112439990
$ws.Cells.Item(145, 3).Value = 17938
$ws.Cells.Item(145, 4).Value = 50166526
$ws.Cells.Item(146, 3).Value = 5160
$ws.Cells.Item(146, 4).Value = 15914236
$ws.Cells.Item(147, 3).Value = 1531
$ws.Cells.Item(147, 4).Value = 6469569
$ws.Cells.Item(148, 3).Value = 344
$ws.Cells.Item(148, 4).Value = 1993100
$ws.Cells.Item(152, 3).Value = 14802
$ws.Cells.Item(152, 4).Value = 22985169
$ws.Cells.Item(153, 3).Value = 53322
$ws.Cells.Item(153, 4).Value = 110424922
$ws.Cells.Item(154, 3).Value = 123864
$ws.Cells.Item(154, 4).Value = 285933329
$ws.Cells.Item(155, 3).Value = 39414
$ws.Cells.Item(155, 4).Value = 113194911
$ws.Cells.Item(156, 3).Value = 11777
$ws.Cells.Item(156, 4).Value = 39622126
$ws.Cells.Item(157, 3).Value = 4218
$ws.Cells.Item(157, 4).Value = 18870920
$ws.Cells.Item(158, 3).Value = 857
$ws.Cells.Item(158, 4).Value = 5222135
$ws.Cells.Item(160, 3).Value = 41085
$ws.Cells.Item(160, 4).Value = 63838834
